# "Doing Updates for Financials" — refresh the SNH yearly financials sheet
# with a new reporting period (inserted as the new column D) and a small
# number of revised prior-period figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Step 1: Insert a new column before D. This shifts the existing D:K
# data right to E:L (a new, currently-empty L column is also created).
# ---------------------------------------------------------------------
$ws.Columns("D:D").Insert()

# ---------------------------------------------------------------------
# Step 2: The freshly inserted column D has no number formatting yet.
# Copy the formatting (date style for the header rows, #,##0 style for
# the data rows) from column E - which still carries the original
# per-row formatting - into the new column D, for each of the three
# contiguous data blocks on the sheet.
# ---------------------------------------------------------------------
$ws.Range("E7:E35").Copy()
$ws.Range("D7:D35").PasteSpecial(-4122)
$ws.Range("E38:E77").Copy()
$ws.Range("D38:D77").PasteSpecial(-4122)
$ws.Range("E80:E102").Copy()
$ws.Range("D80:D102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# Step 3: Populate the new column D with the latest reporting-period
# values (period-ending dates plus Income Statement / Balance Sheet /
# Cash Flow figures).
# ---------------------------------------------------------------------
$ws.Range("D7").Value2 = 43465
$ws.Range("D8").Value2 = 1117200
$ws.Range("D9").Value2 = 451600
$ws.Range("D10").Value2 = 665600
$ws.Range("D12").Value2 = "NA"
$ws.Range("D13").Value2 = 0
$ws.Range("D14").Value2 = 66400
$ws.Range("D15").Value2 = 286200
$ws.Range("D17").Value2 = 890300
$ws.Range("D18").Value2 = 226900
$ws.Range("D20").Value2 = 244800
$ws.Range("D21").Value2 = 757900
$ws.Range("D22").Value2 = 179300
$ws.Range("D23").Value2 = 292400
$ws.Range("D24").Value2 = 500
$ws.Range("D25").Value2 = 0
$ws.Range("D26").Value2 = 291900
$ws.Range("D27").Value2 = 286900
$ws.Range("D28").Value2 = 0
$ws.Range("D29").Value2 = "NA"
$ws.Range("D30").Value2 = 0
$ws.Range("D31").Value2 = 0
$ws.Range("D32").Value2 = -244800
$ws.Range("D33").Value2 = 286900
$ws.Range("D34").Value2 = 0
$ws.Range("D35").Value2 = 286900

$ws.Range("D38").Value2 = 43465
$ws.Range("D41").Value2 = 55000
$ws.Range("D42").Value2 = 0
$ws.Range("D43").Value2 = 18700
$ws.Range("D44").Value2 = 0
$ws.Range("D45").Value2 = 0
$ws.Range("D46").Value2 = 0
$ws.Range("D47").Value2 = 150700
$ws.Range("D48").Value2 = 6341900
$ws.Range("D49").Value2 = 0
$ws.Range("D50").Value2 = 0
$ws.Range("D51").Value2 = 0
$ws.Range("D52").Value2 = 15100
$ws.Range("D53").Value2 = 0
$ws.Range("D54").Value2 = 7160400
$ws.Range("D57").Value2 = 54300
$ws.Range("D58").Value2 = 0
$ws.Range("D59").Value2 = 26200
$ws.Range("D60").Value2 = 0
$ws.Range("D61").Value2 = 3648400
$ws.Range("D62").Value2 = 0
$ws.Range("D63").Value2 = 0
$ws.Range("D64").Value2 = 0
$ws.Range("D65").Value2 = 0
$ws.Range("D66").Value2 = 4137300
$ws.Range("D68").Value2 = 0
$ws.Range("D69").Value2 = 0
$ws.Range("D70").Value2 = 0
$ws.Range("D71").Value2 = 0
$ws.Range("D72").Value2 = -1590400
$ws.Range("D73").Value2 = 0
$ws.Range("D74").Value2 = 0
$ws.Range("D75").Value2 = 0
$ws.Range("D76").Value2 = 3023100
$ws.Range("D77").Value2 = 0

$ws.Range("D80").Value2 = 43465
$ws.Range("D81").Value2 = 286900
$ws.Range("D83").Value2 = 286200
$ws.Range("D84").Value2 = 0
$ws.Range("D85").Value2 = 0
$ws.Range("D86").Value2 = 0
$ws.Range("D87").Value2 = 0
$ws.Range("D88").Value2 = 0
$ws.Range("D89").Value2 = 392800
$ws.Range("D91").Value2 = -233300
$ws.Range("D92").Value2 = 0
$ws.Range("D93").Value2 = 0
$ws.Range("D94").Value2 = 99100
$ws.Range("D96").Value2 = -370700
$ws.Range("D97").Value2 = 0
$ws.Range("D98").Value2 = 0
$ws.Range("D99").Value2 = 0
$ws.Range("D100").Value2 = -469200
$ws.Range("D101").Value2 = 0
$ws.Range("D102").Value2 = 22800

# ---------------------------------------------------------------------
# Step 4: A handful of prior-period figures (now in column E after the
# shift) were also revised as part of this data refresh, not merely
# shifted over from the old column D.
# ---------------------------------------------------------------------
$ws.Range("E8").Value2 = 1074700
$ws.Range("E9").Value2 = 413500
$ws.Range("E10").Value2 = 661200
$ws.Range("E89").Value2 = 419300
$ws.Range("E102").Value2 = 11700
